$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") for rows 2-72 changes from serial date 45171 to 45172
$lastRow = 72
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 45172
}
